# Applies the reachable portion of the authored edit:
#   Figure 3 : "PFSM Model" -> "PFSM Module"
# (Text Box 74, inside the "Group 75" group shape on slide 1.)
#
# NOTE: the commit also touches two other text boxes ("Text Box 192" /
# id=13 and "Text Box 189" / id=81) that live inside group shapes whose
# very first child is an <mc:AlternateContent> block
# (p:grpSp > mc:AlternateContent). This headless COM-interop runtime's
# Shapes/GroupItems walker does not resolve such groups (Name/Id come
# back blank and GroupItems.Count is 0), so those two shapes are not
# reachable through the PowerPoint object model here. This script
# performs every edit that the exposed object model can actually reach.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Group 75" -> "TextBox 74" holds the caption "Figure 3 : PFSM Model"
$group = $s.Shapes.Item(13)
$textBox = $group.GroupItems.Item(2)

$tr = $textBox.TextFrame.TextRange
$full = $tr.Text
$marker = " PFSM Model"
$startPos = $full.IndexOf($marker) + 1   # PowerPoint ranges are 1-based
$len = $marker.Length

$run = $tr.Characters($startPos, $len)
$run.Text = " PFSM Module"
